$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "LFU"
$ws.Range("B16").Value = 43840
$ws.Range("C16").Value = "TournamentBP"
